$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: this engine coalesces adjacent runs that end up with identical
# <w:rPr>, so after doing text edits we re-establish the exact run
# boundaries the target XML needs. Toggling Bold on/off (net no-op, so the
# visible/serialized formatting is unaffected) over a Range forces a new
# run boundary at *both* ends of that Range. Calling it once per segment,
# left to right, with shared endpoints between adjacent segments, produces
# a clean run-per-segment split.
# ---------------------------------------------------------------------------
function Split-Segments($d, $base, [string[]]$segments) {
    $p = 0
    foreach ($seg in $segments) {
        $len = $seg.Length
        $r = $d.Range($base + $p, $base + $p + $len)
        $r.Font.Bold = $true
        $r.Font.Bold = $false
        $p += $len
    }
}

# ===========================================================================
# Change 1: "1-panel advertisements..." -> "Draft copy of 1-panel advertisements..."
#           with the new text split into two runs: "Draft copy of 1" | "-panel..."
# ===========================================================================
$rng = $d.Content
$rng.Find.Execute("1-panel", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$s1 = $rng.Start
$oneChar = $d.Range($s1, $s1 + 1)
$oneChar.Text = "Draft copy of 1"

$rng = $d.Content
$rng.Find.Execute("Draft copy of 1-panel advertisements in multiple media (i.e. page in a magazine, poster on the side of a bus stand, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Split-Segments $d $rng.Start @(
    "Draft copy of 1",
    "-panel advertisements in multiple media (i.e. page in a magazine, poster on the side of a bus stand, "
)

# ===========================================================================
# Change 2: delete the old _GoBack bookmark (currently sits right after the
# "...bus stand, etc)." text in the "Where are we presenting this" paragraph)
# ===========================================================================
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ===========================================================================
# Change 3: "services the most destinations of" -> "operates the most flights of"
#           and re-split the whole "Thanksgiving because it ... 2008." run of
#           text into the correct set of runs (the replace coalesces
#           everything from "Thanksgiving" through "2008." into one run, so
#           we restore every boundary that needs to exist afterwards,
#           including the two brand-new ones).
# ===========================================================================
$rng = $d.Content
$rng.Find.Execute("services the most destinations of", $true, $false, $false, $false, $false, $true, 1, $false, "operates the most flights of", 2)

$rng = $d.Content
$rng.Find.Execute("Thanksgiving because it operates the most flights of all domestic carriers and was one of the leaders in on-time arrivals during Thanksgiving week in 2008", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Split-Segments $d $rng.Start @(
    "Thanksgiving",
    " because it ",
    "operates the most flights of all",
    " domestic carriers and was one of the leaders in on-time arrival",
    "s",
    " during ",
    "Thanksgivi",
    "ng week in 2008."
)

# ===========================================================================
# Change 4: insert a new _GoBack bookmark in the middle of the final
# "Thanksgiving week in 2008." run, between "Thanksgivi" and "ng week in 2008."
# (the run split above already separates these two pieces).
# ===========================================================================
$rng = $d.Content
$rng.Find.Execute("Thanksgiving week in 2008", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmPos = $rng.Start + ("Thanksgivi").Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
